# Apply updated cryptocurrency market data to Sheet1 (rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E (price / volume) columns to plain text so values such as
# "249.07" or "1.00" are stored verbatim instead of being reinterpreted
# as numbers (which would drop trailing zeros / add float noise).
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "37.094.64"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "2.044.68"
$ws.Range("E3").Value = "  -3.42%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "249.07"
$ws.Range("E5").Value = "  -2.73%  "
$ws.Range("D6").Value = "0.657"
$ws.Range("E6").Value = "  -1.81%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "56.07"
$ws.Range("E8").Value = "  +20.58%  "
$ws.Range("D9").Value = "61.90"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "0.379"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").Value = "0.0761"
$ws.Range("E11").Value = "  +2.30%  "
$ws.Range("E12").Value = "  +5.42%  "
$ws.Range("D13").Value = "15.16"
$ws.Range("E13").Value = "  +3.81%  "
$ws.Range("D14").Value = "2.336.81"
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "0.828"
$ws.Range("E15").Value = "  -2.89%  "
$ws.Range("D16").Value = "5.28"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").Value = "2.039.73"
$ws.Range("E17").Value = "  -3.35%  "
$ws.Range("D18").Value = "37.043.33"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "72.59"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").Value = "14.52"
$ws.Range("E20").Value = "  +7.62%  "
$ws.Range("D21").Value = "0.0₃0867"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("D22").Value = "238.58"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("D23").Value = "5.25"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "2.42"
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("D26").Value = "170.33"
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("D27").Value = "9.14"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").Value = "20.36"
$ws.Range("E28").Value = "  -4.86%  "
$ws.Range("D29").Value = "2.01"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").Value = "1.08"
$ws.Range("E31").Value = "  +17.90%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.57"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("B33").Value = "Gas"
$ws.Range("C33").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D33").Value = "22.10"
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("D34").Value = "0.0632"
$ws.Range("E34").Value = "  +4.48%  "
$ws.Range("E35").Value = "  +4.15%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "2.30"
$ws.Range("E37").Value = "  -4.80%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.0862"
$ws.Range("E38").Value = "  -9.80%  "
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  -7.12%  "
$ws.Range("D40").Value = "1.36"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").Value = "0.107"
$ws.Range("E41").Value = "  +28.12%  "
$ws.Range("D42").Value = "18.31"
$ws.Range("E42").Value = "  +12.42%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "4.61"
$ws.Range("E43").Value = "  +81.36%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0226"
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "1.15"
$ws.Range("E45").Value = "  -4.49%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "98.14"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D47").Value = "2.80"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").Value = "1.306.35"
$ws.Range("E48").Value = "  -4.35%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "2.39"
$ws.Range("E49").Value = "  +3.19%  "
$ws.Range("D50").Value = "2.90"
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("D51").Value = "6.91"
$ws.Range("E51").Value = "  -0.18%  "

# Restore the default (unstyled) look now that the literal text is in place.
$dataRange.Style = "Normal"
